$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Rename test case identifiers to match the "TC_SM_COA_*" naming convention
# used elsewhere in the sheet (accessibility/Axe pass over the MDOT screens).
$ws.Range("B9").Value  = "TC_SM_COA_AppropriationNumberProfile_EditProfile"
$ws.Range("B15").Value = "TC_SM_COA_DeleteProfile"

# The delete-profile step now runs for 2 iterations instead of 1.
$ws.Range("F15").Value = 2

# Leave the selection where the author left it when saving.
$ws.Range("E19").Select()
